$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")
$ws.Activate()

# --- Append ORDER BY / LIMIT clauses to the Cypher queries in column B ---
# CasesTab query (B2)
$ws.Range("B2").Value2 = $ws.Range("B2").Value2 + "`n order By ss.study_subject_id ASC LIMIT 100"

# SamplesTab query (B3)
$ws.Range("B3").Value2 = $ws.Range("B3").Value2 + "`n order By samp.sample_id ASC LIMIT 100"

# FilesTab query (B4) - replace the existing lowercase "order by" clause
$b4 = $ws.Range("B4").Value2
$b4 = $b4.Replace("    order by f.file_name", "     order By f.file_name ASC LIMIT 100")
$ws.Range("B4").Value2 = $b4

# --- Row heights grew slightly because the wrapped text got longer ---
$ws.Rows.Item(2).RowHeight = 360
$ws.Rows.Item(3).RowHeight = 360

# --- Update the window/sheet view state (scroll + selection) ---
$ws.Range("A3").Select()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1

$ws.Range("B4").Select()

$wb.Save()
